$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

# Drop the "[Dev] " prefix that used to mark these as dev-only asset names.
$settings.Range("B2").Value = "RPA_Moon_UploadBucket"

$settings.Range("B7").Value = "RPA_Moon_Cred_Gmail"
$settings.Range("B7").Style = "Normal"

$settings.Range("B8").Value = "RPA090_Expedia_FlagError"
$settings.Range("B8").Style = "Normal"

# ---------------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")

# Selection on Settings moved from B8 to A8:B8 (A8 is now the active cell).
# Activate Settings only long enough to move its selection, then restore
# Assets as the active tab (it is the tab that was active originally).
$settings.Activate()
$settings.Range("A8:B8").Select()
$assets.Activate()

$assets.Range("B2").Value = "RPA_Moon_SheetIdConfig_Accommodation"
$assets.Range("B3").Value = "RPA_Moon_PathMasterFolder"
# B4 keeps its pre-existing quote-prefix ('text) cell style.
$assets.Range("B4").Value = "'RPA_Moon_PathMailTemplate"
$assets.Range("B5").Value = "RPA_Moon_PathSaKey"

# Row 6/7/8 get repurposed for the new Captcha-related assets, and the old
# PathDownloadChrome / DialogDownloadChrome / SenderName rows move down to
# rows 9/10/11.
$assets.Range("A6").Value = "Captcha_SiteKey"
$assets.Range("B6").Value = "RPA021_MOLPAY_Captcha_SiteKey"

$assets.Range("A7").Value = "Captcha_RuleId"
$assets.Range("B7").Value = "RPA_Moon_Captcha_RuleId"

$assets.Range("A8").Value = "GCaptcha_RuleId"
$assets.Range("B8").Value = "RPA_Moon_GCaptcha_RuleId"

$assets.Range("A9").Value = "PathDownloadChrome"
$assets.Range("B9").Value = "RPA_Moon_PathDownloadChrome"

$assets.Range("A10").Value = "DialogDownloadChrome"
$assets.Range("B10").Value = "RPA_Moon_DialogDownloadChrome"

$assets.Range("A11").Value = "SenderName"
$assets.Range("B11").Value = "RPA_Moon_SenderName"

# Extend the sheet with a few more blank formatted rows at the bottom.
$assets.Rows.Item(998).RowHeight = 14.3
$assets.Rows.Item(999).RowHeight = 14.3
$assets.Rows.Item(1000).RowHeight = 14.3
